$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) First paragraph: "This is a Microsoft word document." gains a
#    trailing double-space and a new red parenthetical remark, typed
#    as three consecutive runs (matching how the edit was captured).
# ------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$insertPoint = $p1.End - 1
$spaceRng = $d.Range($insertPoint, $insertPoint)
$spaceRng.InsertAfter("  ")

$ip2 = $d.Paragraphs(1).Range.End - 1
$runA = $d.Range($ip2, $ip2)
$runA.InsertAfter("(This is a change – Ve")
$runA.Font.Color = 255

$ip3 = $d.Paragraphs(1).Range.End - 1
$runB = $d.Range($ip3, $ip3)
$runB.InsertAfter("rsion for main branch")
$runB.Font.Color = 255

$ip4 = $d.Paragraphs(1).Range.End - 1
$runC = $d.Range($ip4, $ip4)
$runC.InsertAfter(")")
$runC.Font.Color = 255

# ------------------------------------------------------------------
# 2) Drop the trailing "ank God almighty, we are free at last."
#    paragraph (the stray leftover of a pasted web snippet) that sat
#    just after the poem's last line.
# ------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$d.Paragraphs($lastIndex).Range.Delete()

# ------------------------------------------------------------------
# 3) That pasted snippet is also the only thing in the document that
#    used the orphaned styles below (plus a couple of unrelated
#    leftover Heading styles) - remove them now that nothing
#    references them. Deleting from the highest style index down
#    keeps earlier indices stable while we work.
# ------------------------------------------------------------------
$staleStyleIndexes = @(18, 17, 16, 15, 14, 13, 12, 11, 10, 3, 2)
foreach ($idx in $staleStyleIndexes) {
    $d.Styles($idx).Delete()
}
